$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Banda in download (Mb/s)"
$ws.Range("C1").Value = "Banda in upload (Mb/s)"
$ws.Range("D1").Value = "Tempo totale download (s)"
$ws.Range("E1").Value = "Tempo totale upload (s)"

$ws.Range("D2").Value = 16.002
$ws.Range("E2").Value = 22.028
$ws.Range("D3").Value = 19.006
$ws.Range("E3").Value = 25.027
$ws.Range("D4").Value = 19.002
$ws.Range("E4").Value = 31.038
